$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 3 brand-new rows at the top of the "Alcachofa" block (rows
#    242-244), pushing the existing rows 242-272 down to 245-275.
# ---------------------------------------------------------------------------
$ws.Range("A242:R244").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2) Populate the 3 newly-inserted rows with the new weekly data
#    (Fecha = 2021-09-09 -> serial 44449).
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=242; A=9; B="Vega Central Mapocho de Santiago"; C="Metropolitana"; D=44449; E=13; F=100112013; G="Alcachofa"; H="Española"; I="Extra";   J=16; K=14000; L=15000; M=14500; N="`$/caja 25 unidades"; O="Provincia de Limarí"; P=14500; Q=1;  R="Hortaliza" },
    @{ Row=243; A=9; B="Vega Central Mapocho de Santiago"; C="Metropolitana"; D=44449; E=13; F=100112013; G="Alcachofa"; H="Española"; I="Primera"; J=43; K=12000; L=13000; M=12488; N="`$/caja 30 unidades"; O="Provincia de Limarí"; P=416;   Q=30; R="Hortaliza" },
    @{ Row=244; A=9; B="Vega Central Mapocho de Santiago"; C="Metropolitana"; D=44449; E=13; F=100112013; G="Alcachofa"; H="Española"; I="Segunda"; J=25; K=10000; L=11000; M=10480; N="`$/caja 40 unidades"; O="Provincia de Limarí"; P=262;   Q=40; R="Hortaliza" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $r.A
    $ws.Cells.Item($row, 2).Value2  = $r.B
    $ws.Cells.Item($row, 3).Value2  = $r.C
    $ws.Cells.Item($row, 4).Value2  = $r.D
    $ws.Cells.Item($row, 5).Value2  = $r.E
    $ws.Cells.Item($row, 6).Value2  = $r.F
    $ws.Cells.Item($row, 7).Value2  = $r.G
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
}

# ---------------------------------------------------------------------------
# 3) The "Española" block that used to sit at old rows 245-246 (now shifted
#    to rows 248-249) had its Fecha corrected from 2021-08-30 (44428) to
#    2021-09-09 (44438).
# ---------------------------------------------------------------------------
$ws.Range("D248").Value2 = 44438
$ws.Range("D249").Value2 = 44438
